$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.628542900085449
$ws.Range("B1").Value = 1.820488452911377
$ws.Range("C1").Value = 4.984848976135254
$ws.Range("D1").Value = 1.453904151916504
$ws.Range("E1").Value = 0.7635625600814819
